# Updated Audio asset List
# Applies the edits described by the commit to the "Stuart Audio Asset List" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F/G "Wwise/Unity STATUS" cells flipped from Not implemented -> Complete ---
$ws.Range("F14").Value = "Complete"
$ws.Range("F14").Style = "Good"

$ws.Range("F16").Value = "Complete"
$ws.Range("F16").Style = "Good"
$ws.Range("G16").Value = "Complete"
$ws.Range("G16").Style = "Good"

$ws.Range("F17").Value = "Complete"
$ws.Range("F17").Style = "Good"
$ws.Range("G17").Value = "Complete"
$ws.Range("G17").Style = "Good"

$ws.Range("F19").Value = "Complete"
$ws.Range("F19").Style = "Good"

$ws.Range("F24").Value = "Complete"
$ws.Range("F24").Style = "Good"
$ws.Range("G24").Value = "Complete"
$ws.Range("G24").Style = "Good"

$ws.Range("F25").Value = "Complete"
$ws.Range("F25").Style = "Good"

# --- Rows 27-35: level-theme rows, now resolved (row heights shrink, description updated) ---

# Row 27: theme decided -> Aztec/Mayan Tropes
$ws.Range("C27").Value = "Aztec/Mayan Tropes"
$ws.Range("C27").Style = "Normal"
$ws.Rows.Item(27).RowHeight = 30

# Row 28: theme decided -> Egyptian/Middle Eastern Tropes
$ws.Range("C28").Value = "Egyptian/Middle Eastern Tropes"
$ws.Range("C28").Style = "Normal"
$ws.Range("E28").Value = "Complete"
$ws.Range("E28").Style = "Good"
$ws.Range("F28").Value = "Complete"
$ws.Range("F28").Style = "Good"
$ws.Range("G28").Value = "Complete"
$ws.Range("G28").Style = "Good"
$ws.Rows.Item(28).RowHeight = 30

# Row 29: theme decided -> Far Eastern/Asian Tropes
$ws.Range("C29").Value = "Far Eastern/Asian Tropes"
$ws.Range("C29").Style = "Normal"
$ws.Rows.Item(29).RowHeight = 30

# Row 30: theme still awaiting confirmation -> clear description text, keep style, progress updated
$ws.Range("C30").ClearContents()
$ws.Range("E30").Value = "Complete"
$ws.Range("E30").Style = "Good"
$ws.Range("F30").Value = "Complete"
$ws.Range("F30").Style = "Good"
$ws.Rows.Item(30).RowHeight = 30

# Row 31
$ws.Range("C31").ClearContents()
$ws.Range("E31").Value = "Complete"
$ws.Range("E31").Style = "Good"
$ws.Range("F31").Value = "Complete"
$ws.Range("F31").Style = "Good"
$ws.Rows.Item(31).RowHeight = 30

# Row 32
$ws.Range("C32").ClearContents()
$ws.Range("E32").Value = "Complete"
$ws.Range("E32").Style = "Good"
$ws.Range("F32").Value = "Complete"
$ws.Range("F32").Style = "Good"
$ws.Rows.Item(32).RowHeight = 30

# Row 33
$ws.Range("C33").ClearContents()
$ws.Range("E33").Value = "Complete"
$ws.Range("E33").Style = "Good"
$ws.Range("F33").Value = "Complete"
$ws.Range("F33").Style = "Good"
$ws.Rows.Item(33).RowHeight = 30

# Row 34
$ws.Range("C34").ClearContents()
$ws.Range("E34").Value = "Complete"
$ws.Range("E34").Style = "Good"
$ws.Range("F34").Value = "Complete"
$ws.Range("F34").Style = "Good"
$ws.Rows.Item(34).RowHeight = 30

# Row 35
$ws.Range("C35").ClearContents()
$ws.Range("E35").Value = "Complete"
$ws.Range("E35").Style = "Good"
$ws.Range("F35").Value = "Complete"
$ws.Range("F35").Style = "Good"
$ws.Rows.Item(35).RowHeight = 30

# --- Update the active selection to match the saved view state ---
$ws.Range("F14").Select()
